$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")

$ws.Range("D2").Value = 4430.59
$ws.Range("E2").Value = -4430.59

$ws.Range("D3").Value = 594.78
$ws.Range("E3").Value = 16905.22
$ws.Range("F3").Value = 0.03398742857142857

$ws.Range("D4").Value = 5025.37
$ws.Range("E4").Value = 12474.63
$ws.Range("F4").Value = 0.287164
